$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture a pristine default (unstyled) cell style so we can restore it
# after forcing the Price column to Text format (prevents Excel from
# auto-converting decimal-looking strings like "31.96" into numbers).
$defaultStyle = $ws.Range("Z1").Style
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "34.115.14"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.786.96"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "226.72"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "31.96"
$ws.Range("E8").Value = "  -2.17%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("D12").Value = "2.044.36"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "11.25"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "1.794.55"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "34.029.95"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "68.01"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "245.96"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "4.10"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  -2.04%  "
$ws.Range("D25").Value = "161.66"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "7.16"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("E33").Value = "  +2.39%  "
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "1.447.89"
$ws.Range("D36").Value = "0.647"
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("D38").Value = "2.40"
$ws.Range("E38").Value = "  +7.80%  "
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "80.13"
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").Value = "13.56"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("E45").Value = "  +1.93%  "
$ws.Range("D46").Value = "6.04"
$ws.Range("E46").Value = "  +3.38%  "
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").Value = "107.55"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("D50").Value = "1.945.49"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  +0.26%  "

# Restore original (default) styling on the Price column
$ws.Range("D2:D51").Style = $defaultStyle
